$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 22:50"

# Swap the order of Noruega/Brasil: Brasil's row now comes first (row 21) with
# new case numbers, and Noruega (the old row-21 data) moves down to row 22.
$ws.Range("A21").Value = "Brasil"
$ws.Range("B21").Value = 4579
$ws.Range("C21").Value = 323
$ws.Range("D21").Value = 120
$ws.Range("E21").Value = 4300
$ws.Range("F21").Value = 296
$ws.Range("G21").Value = 23
$ws.Range("H21").Value = 159

$ws.Range("A22").Value = "Noruega"
$ws.Range("B22").Value = 4445
$ws.Range("C22").Value = 161
$ws.Range("D22").Value = 12
$ws.Range("E22").Value = 4401
$ws.Range("F22").Value = 97
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 32

# Chequia (row 25) case-count update
$ws.Range("B25").Value = 2976
$ws.Range("C25").Value = 159
$ws.Range("E25").Value = 2928
